$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ducry1979")

# Update the upper bound of each 10-year age bracket (G column) from X9 to X0
# for the two blocks of rows (37-44 and 71-78), e.g. 20-29 -> 20-30
foreach ($row in 37..44) {
    $c = $ws.Cells.Item($row, 7)
    $c.Value2 = $c.Value2 + 1
}
foreach ($row in 71..78) {
    $c = $ws.Cells.Item($row, 7)
    $c.Value2 = $c.Value2 + 1
}
